# Update results for DA with asset smoothing
# (Conditional indexation column, L, on both result sheets)

$wb = $excel.ActiveWorkbook

$ws15 = $wb.Worksheets.Item("ERCvol_15y")
$ws15.Range("L2").Value = 9.68125932284621
$ws15.Range("L3").Value = 6.830689243486568
$ws15.Range("L4").Value = 4.678333523980239
$ws15.Range("L5").Value = 1.6009054314387683

$ws30 = $wb.Worksheets.Item("ERCvol_30y")
$ws30.Range("L2").Value = 13.107931845338065
$ws30.Range("L3").Value = 12.583896662614858
$ws30.Range("L4").Value = 10.712680603328597
$ws30.Range("L5").Value = 7.7846627391028616
$ws30.Range("L6").Value = 1.783263634148224
